$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new "time_taken" column header in F1, matching the header style used
# by the existing header cells (B1:E1) -- copy E1's formatting onto F1
# (reuses the existing header cell style rather than minting a new one).
$ws.Range("F1").Value = "time_taken"
$ws.Range("E1").Copy()
$ws.Range("F1").PasteSpecial(-4122)

# Populate F2:F12 with the per-row "time_taken" metadata timestamps.
$ws.Range("F2").Value = "2021-10-05 13:38:37.850234"
$ws.Range("F3").Value = "2021-10-05 13:38:37.850247"
$ws.Range("F4").Value = "2021-10-05 13:38:37.850251"
$ws.Range("F5").Value = "2021-10-05 13:38:37.850254"
$ws.Range("F6").Value = "2021-10-05 13:38:37.850257"
$ws.Range("F7").Value = "2021-10-05 13:38:37.850260"
$ws.Range("F8").Value = "2021-10-05 13:38:37.850263"
$ws.Range("F9").Value = "2021-10-05 13:38:37.850266"
$ws.Range("F10").Value = "2021-10-05 13:38:37.850269"
$ws.Range("F11").Value = "2021-10-05 13:38:37.850273"
$ws.Range("F12").Value = "2021-10-05 13:38:37.850275"
